# Add files via upload
# Adds a new "experimental_results" worksheet (absorbance / concentration data
# for the LARGE and SMALL particle runs) as the last sheet in the workbook,
# and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet after the current last sheet (T_test) and
#    rename it.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "experimental_results"

# ---------------------------------------------------------------------------
# 2. Column widths for the descriptive / concentration columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 46.83
$ws.Columns.Item(4).ColumnWidth = 11.67
$ws.Columns.Item(7).ColumnWidth = 46.83
$ws.Columns.Item(8).ColumnWidth = 11.67

# ---------------------------------------------------------------------------
# 3. Row 2 header labels (Time / Absorbance / Concentration (.../ Concentration)
#    for both the LARGE (A:D) and SMALL (E:H) blocks.
# ---------------------------------------------------------------------------
$headers = New-Object 'object[,]' 1,4
$headers[0,0] = "Time"
$headers[0,1] = "Absorbance"
$headers[0,2] = "Concentration (measured using standard curve, diluted)"
$headers[0,3] = "concentration"

$ws.Range("A2:D2").Value = $headers
$ws.Range("E2:H2").Value = $headers

# ---------------------------------------------------------------------------
# 4. Data rows 3-11: LARGE block in A:D, SMALL block in E:H.
# ---------------------------------------------------------------------------
$large = New-Object 'object[,]' 9,4
$large[0,0]=30;  $large[0,1]=1.8029999999999999;  $large[0,2]=0.42366336630000001; $large[0,3]=0.84732673270000003
$large[1,0]=60;  $large[1,1]=1.44;                 $large[1,2]=0.33381188119999999; $large[1,3]=0.66762376239999999
$large[2,0]=90;  $large[2,1]=1.7669999999999999;  $large[2,2]=0.41475247520000003; $large[2,3]=0.82950495049999995
$large[3,0]=30;  $large[3,1]=1.4950000000000001;  $large[3,2]=0.34742574259999998; $large[3,3]=0.69485148510000005
$large[4,0]=60;  $large[4,1]=1.506;                 $large[4,2]=0.35014851489999999; $large[4,3]=0.70029702969999996
$large[5,0]=90;  $large[5,1]=1.63;                  $large[5,2]=0.38084158420000003; $large[5,3]=0.76168316830000005
$large[6,0]=30;  $large[6,1]=1.389;                 $large[6,2]=0.32118811879999998; $large[6,3]=0.64237623759999996
$large[7,0]=60;  $large[7,1]=2.4049999999999998;  $large[7,2]=0.57267326730000001; $large[7,3]=1.1453465350000001
$large[8,0]=90;  $large[8,1]=2.6219999999999999;  $large[8,2]=0.62638613860000003; $large[8,3]=1.252772277

$ws.Range("A3:D11").Value = $large

$small = New-Object 'object[,]' 9,4
$small[0,0]=30;  $small[0,1]=1.0349999999999999;  $small[0,2]=0.23356435640000001; $small[0,3]=0.46712871290000002
$small[1,0]=60;  $small[1,1]=1.45;                  $small[1,2]=0.33628712869999999; $small[1,3]=0.67257425739999999
$small[2,0]=90;  $small[2,1]=0.97199999999999998;  $small[2,2]=0.21797029700000001; $small[2,3]=0.43594059410000002
$small[3,0]=30;  $small[3,1]=1.0569999999999999;  $small[3,2]=0.239009901;          $small[3,3]=0.47801980199999999
$small[4,0]=60;  $small[4,1]=1.21;                  $small[4,2]=0.27688118810000001; $small[4,3]=0.55376237620000002
$small[5,0]=90;  $small[5,1]=0.79;                  $small[5,2]=0.1729207921;         $small[5,3]=0.3458415842
$small[6,0]=30;  $small[6,1]=0.88100000000000001;  $small[6,2]=0.19544554459999999; $small[6,3]=0.39089108909999998
$small[7,0]=60;  $small[7,1]=1.1990000000000001;  $small[7,2]=0.2741584158;         $small[7,3]=0.54831683170000001
$small[8,0]=90;  $small[8,1]=0.76700000000000002;  $small[8,2]=0.16722772280000001; $small[8,3]=0.3344554455

$ws.Range("E3:H11").Value = $small

# ---------------------------------------------------------------------------
# 5. Row 1 merged title cells ("LARGE" over A:D, "SMALL" over E:H).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "LARGE"
$ws.Range("E1").Value = "SMALL"
$ws.Range("A1:D1").Merge()
$ws.Range("E1:H1").Merge()

# ---------------------------------------------------------------------------
# 6. Styling.
#    - Data rows (3-11) reuse the workbook's existing "Arial 10" style (the
#      same style already used on the small/large/T_test sheets).
#    - Row 2 headers are bold Arial 10.
#    - Row 1 merged titles are bold Calibri 12, centered.
#    A couple of helper ("seed") cells build each style cleanly (matching an
#    existing cell's format, then tweaking just one attribute) so the
#    workbook doesn't accumulate duplicate/intermediate styles; they are
#    cleared again once used.
# ---------------------------------------------------------------------------

# Data-row style: copy the existing Arial 10 style used on the "small" sheet.
$wb.Worksheets.Item("small").Range("A2").Copy() | Out-Null
$ws.Range("A3:H11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Seed Z1: bold Arial 10 (copy existing Arial 10 style, then bold it).
$wb.Worksheets.Item("small").Range("A2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("Z1").Font.Bold = $true

# Seed Z2: bold Calibri 12, centered (default font + bold + center).
$ws.Range("Z2").Font.Bold = $true
$ws.Range("Z2").HorizontalAlignment = -4108  # xlCenter

# Apply the seeded styles to the real ranges, then clean up the seed cells.
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("A2:H2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("Z2").Copy() | Out-Null
$ws.Range("A1:D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("E1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("Z1:Z2").Clear() | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 7. Final selection / activation state.
# ---------------------------------------------------------------------------
$ws.Range("G13").Select()
$ws.Activate()
